# Updates the pl_mw results table (rows 2-25, columns B:O) with the
# recomputed values for the "case with 380 kV" run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nRows = 24
$nCols = 14
$arr = New-Object 'object[,]' $nRows,$nCols

$arr[0,0] = 0.488581626683299
$arr[0,1] = 0.09635221466777466
$arr[0,2] = 0.07767856683607022
$arr[0,3] = 0.1046558067673296
$arr[0,4] = 0
$arr[0,5] = 0.002502358733471786
$arr[0,6] = 0
$arr[0,7] = 1.229678219994042
$arr[0,8] = 0
$arr[0,9] = 0.3093684014095288
$arr[0,10] = 0.2166115354202276
$arr[0,11] = 0.1549274488537939
$arr[0,12] = 0
$arr[0,13] = 4.971152512444831

$arr[1,0] = 0.4561002489688519
$arr[1,1] = 0.09359144670311537
$arr[1,2] = 0.07058895358652251
$arr[1,3] = 0.1049968585447729
$arr[1,4] = 0
$arr[1,5] = 0.002505151990921543
$arr[1,6] = 0
$arr[1,7] = 1.235913461079477
$arr[1,8] = 0
$arr[1,9] = 0.2769684345009011
$arr[1,10] = 0.213987690322476
$arr[1,11] = 0.148641733463208
$arr[1,12] = 0
$arr[1,13] = 4.985913092587481

$arr[2,0] = 0.4363306699750069
$arr[2,1] = 0.09187076005728301
$arr[2,2] = 0.06627036952073695
$arr[2,3] = 0.1052417477820153
$arr[2,4] = 0
$arr[2,5] = 0.002506959422347355
$arr[2,6] = 0
$arr[2,7] = 1.240257398778358
$arr[2,8] = 0
$arr[2,9] = 0.2571214955987386
$arr[2,10] = 0.2124806973957618
$arr[2,11] = 0.1448520100549899
$arr[2,12] = 0
$arr[2,13] = 4.997261807856091

$arr[3,0] = 0.4283186876235163
$arr[3,1] = 0.09116314243374291
$arr[3,2] = 0.0645191829371754
$arr[3,3] = 0.105350480245459
$arr[3,4] = 0
$arr[3,5] = 0.002507719259388627
$arr[3,6] = 0
$arr[3,7] = 1.242157258120027
$arr[3,8] = 0
$arr[3,9] = 0.2490458493446965
$arr[3,10] = 0.2118928096694077
$arr[3,11] = 0.1433253028421042
$arr[3,12] = 0
$arr[3,13] = 5.002461410806632

$arr[4,0] = 0.426990993268447
$arr[4,1] = 0.09104525521696871
$arr[4,2] = 0.06422892361102583
$arr[4,3] = 0.1053690754944405
$arr[4,4] = 0
$arr[4,5] = 0.002507846838535132
$arr[4,6] = 0
$arr[4,7] = 1.242480561489227
$arr[4,8] = 0
$arr[4,9] = 0.2477056390859644
$arr[4,10] = 0.2117967769356923
$arr[4,11] = 0.1430728624256048
$arr[4,12] = 0
$arr[4,13] = 5.00335952883654

$arr[5,0] = 0.436222437621808
$arr[5,1] = 0.091861242882878
$arr[5,2] = 0.06624671729001363
$arr[5,3] = 0.1052431779747796
$arr[5,4] = 0
$arr[5,5] = 0.002506969575313105
$arr[5,6] = 0
$arr[5,7] = 1.240282495855844
$arr[5,8] = 0
$arr[5,9] = 0.2570125349372034
$arr[5,10] = 0.2124726626786568
$arr[5,11] = 0.1448313488058979
$arr[5,12] = 0
$arr[5,13] = 4.997329603528215

$arr[6,0] = 0.4773461968112258
$arr[6,1] = 0.0954056095900242
$arr[6,2] = 0.07522691502248335
$arr[6,3] = 0.1047660465379021
$arr[6,4] = 0
$arr[6,5] = 0.002503302721253626
$arr[6,6] = 0
$arr[6,7] = 1.231721173026308
$arr[6,8] = 0
$arr[6,9] = 0.2981874034840928
$arr[6,10] = 0.2156852695682048
$arr[6,11] = 0.1527457198982134
$arr[6,12] = 0
$arr[6,13] = 4.975767661212217

$arr[7,0] = 0.5593538685816952
$arr[7,1] = 0.1021533557518666
$arr[7,2] = 0.09311128709703098
$arr[7,3] = 0.1041112720091206
$arr[7,4] = 0
$arr[7,5] = 0.00249684169707835
$arr[7,6] = 0
$arr[7,7] = 1.219020834769914
$arr[7,8] = 0
$arr[7,9] = 0.3792896902130849
$arr[7,10] = 0.2228088706492883
$arr[7,11] = 0.168815552137584
$arr[7,12] = 0
$arr[7,13] = 4.951617123272769

$arr[8,0] = 0.6204196022040662
$arr[8,1] = 0.1069878730216516
$arr[8,2] = 0.1064205225980572
$arr[8,3] = 0.1038006199197632
$arr[8,4] = 0
$arr[8,5] = 0.00249253516640444
$arr[8,6] = 0
$arr[8,7] = 1.212180779913702
$arr[8,8] = 0
$arr[8,9] = 0.4390832307544201
$arr[8,10] = 0.2285429154188989
$arr[8,11] = 0.1809538843428911
$arr[8,12] = 0
$arr[8,13] = 4.944928870895666

$arr[9,0] = 0.6483737320904197
$arr[9,1] = 0.1091606435780363
$arr[9,2] = 0.1125126523666182
$arr[9,3] = 0.1036961301783119
$arr[9,4] = 0
$arr[9,5] = 0.00249067069495744
$arr[9,6] = 0
$arr[9,7] = 1.209609680974708
$arr[9,8] = 0
$arr[9,9] = 0.4663280485778216
$arr[9,10] = 0.2312597942434991
$arr[9,11] = 0.1865473618504652
$arr[9,12] = 0
$arr[9,13] = 4.944287417110019

$arr[10,0] = 0.6589839929633854
$arr[10,1] = 0.1099796098760351
$arr[10,2] = 0.1148250184733683
$arr[10,3] = 0.103661843551734
$arr[10,4] = 0
$arr[10,5] = 0.002489978198920868
$arr[10,6] = 0
$arr[10,7] = 1.208713764195743
$arr[10,8] = 0
$arr[10,9] = 0.4766510446440577
$arr[10,10] = 0.2323041526878029
$arr[10,11] = 0.1886756929216276
$arr[10,12] = 0
$arr[10,13] = 4.944389752327908

$arr[11,0] = 0.656697796044341
$arr[11,1] = 0.1098034006122361
$arr[11,2] = 0.1143267686888123
$arr[11,3] = 0.1036689931284887
$arr[11,4] = 0
$arr[11,5] = 0.002490126739246719
$arr[11,6] = 0
$arr[11,7] = 1.208903260333045
$arr[11,8] = 0
$arr[11,9] = 0.4744275435992904
$arr[11,10] = 0.2320785413804316
$arr[11,11] = 0.1882168666745869
$arr[11,12] = 0
$arr[11,13] = 4.944352358457223

$arr[12,0] = 0.6492461538696546
$arr[12,1] = 0.1092280969247241
$arr[12,2] = 0.1127027837844281
$arr[12,3] = 0.103693203632611
$arr[12,4] = 0
$arr[12,5] = 0.002490613452092671
$arr[12,6] = 0
$arr[12,7] = 1.209534416184844
$arr[12,8] = 0
$arr[12,9] = 0.467177210304385
$arr[12,10] = 0.2313454032161673
$arr[12,11] = 0.1867222570627689
$arr[12,12] = 0
$arr[12,13] = 4.944288917344863

$arr[13,0] = 0.6446850008670708
$arr[13,1] = 0.1088752097853103
$arr[13,2] = 0.1117087501775842
$arr[13,3] = 0.1037087206297151
$arr[13,4] = 0
$arr[13,5] = 0.002490913338989115
$arr[13,6] = 0
$arr[13,7] = 1.209931136029688
$arr[13,8] = 0
$arr[13,9] = 0.4627369376198658
$arr[13,10] = 0.2308983563921885
$arr[13,11] = 0.1858080910394762
$arr[13,12] = 0
$arr[13,13] = 4.944295017342824

$arr[14,0] = 0.6185962106653449
$arr[14,1] = 0.1068453439483648
$arr[14,2] = 0.1060231450987033
$arr[14,3] = 0.1038081884035211
$arr[14,4] = 0
$arr[14,5] = 0.002492658911939582
$arr[14,6] = 0
$arr[14,7] = 1.212359681562894
$arr[14,8] = 0
$arr[14,9] = 0.4373035764968449
$arr[14,10] = 0.2283675380646173
$arr[14,11] = 0.1805897709131585
$arr[14,12] = 0
$arr[14,13] = 4.945019104385835

$arr[15,0] = 0.6026360328945088
$arr[15,1] = 0.1055933005380325
$arr[15,2] = 0.102544856508203
$arr[15,3] = 0.103878630821022
$arr[15,4] = 0
$arr[15,5] = 0.002493753946444422
$arr[15,6] = 0
$arr[15,7] = 1.213987932215531
$arr[15,8] = 0
$arr[15,9] = 0.4217121186783572
$arr[15,10] = 0.2268426946583872
$arr[15,11] = 0.177406786180633
$arr[15,12] = 0
$arr[15,13] = 4.946078257226361

$arr[16,0] = 0.5934726661378704
$arr[16,1] = 0.1048706667243522
$arr[16,2] = 0.1005477853531005
$arr[16,3] = 0.1039226146686509
$arr[16,4] = 0
$arr[16,5] = 0.002494392688298611
$arr[16,6] = 0
$arr[16,7] = 1.214975332127366
$arr[16,8] = 0
$arr[16,9] = 0.4127485317976038
$arr[16,10] = 0.2259758547181008
$arr[16,11] = 0.1755827716269351
$arr[16,12] = 0
$arr[16,13] = 4.946913455395247

$arr[17,0] = 0.5903729585726865
$arr[17,1] = 0.1046255676371146
$arr[17,2] = 0.09987222115461236
$arr[17,3] = 0.1039381027824184
$arr[17,4] = 0
$arr[17,5] = 0.002494610487192784
$arr[17,6] = 0
$arr[17,7] = 1.215318386824691
$arr[17,8] = 0
$arr[17,9] = 0.4097143505690042
$arr[17,10] = 0.2256841128676683
$arr[17,11] = 0.1749663549017981
$arr[17,12] = 0
$arr[17,13] = 4.94723505605657

$arr[18,0] = 0.604333316662121
$arr[18,1] = 0.1057268406201501
$arr[18,2] = 0.1029147590960093
$arr[18,3] = 0.1038707733641822
$arr[18,4] = 0
$arr[18,5] = 0.002493636456685169
$arr[18,6] = 0
$arr[18,7] = 1.213809337091547
$arr[18,8] = 0
$arr[18,9] = 0.4233714235542436
$arr[18,10] = 0.2270039605951126
$arr[18,11] = 0.1777449220885288
$arr[18,12] = 0
$arr[18,13] = 4.9459421185054

$arr[19,0] = 0.6514342173326781
$arr[19,1] = 0.1093971811535397
$arr[19,2] = 0.11317964081627
$arr[19,3] = 0.1036859492002513
$arr[19,4] = 0
$arr[19,5] = 0.002490470125637119
$arr[19,6] = 0
$arr[19,7] = 1.20934692172235
$arr[19,8] = 0
$arr[19,9] = 0.4693066500892371
$arr[19,10] = 0.2315603224962359
$arr[19,11] = 0.1871609840879032
$arr[19,12] = 0
$arr[19,13] = 4.944298182070554

$arr[20,0] = 0.6823607467976274
$arr[20,1] = 0.1117737220925221
$arr[20,2] = 0.1199198552452003
$arr[20,3] = 0.1035959331432466
$arr[20,4] = 0
$arr[20,5] = 0.002488479630645342
$arr[20,6] = 0
$arr[20,7] = 1.206883377590906
$arr[20,8] = 0
$arr[20,9] = 0.4993626104191549
$arr[20,10] = 0.2346286940314997
$arr[20,11] = 0.1933743439263154
$arr[20,12] = 0
$arr[20,13] = 4.945236082412066

$arr[21,0] = 0.6658417322196328
$arr[21,1] = 0.1105073548055344
$arr[21,2] = 0.1163195956806362
$arr[21,3] = 0.1036411650314637
$arr[21,4] = 0
$arr[21,5] = 0.002489534798048216
$arr[21,6] = 0
$arr[21,7] = 1.208156782525961
$arr[21,8] = 0
$arr[21,9] = 0.4833181464590268
$arr[21,10] = 0.2329827821852888
$arr[21,11] = 0.1900527550679101
$arr[21,12] = 0
$arr[21,13] = 4.944551391152828

$arr[22,0] = 0.6035659361568548
$arr[22,1] = 0.1056664759250765
$arr[22,2] = 0.1027475179029267
$arr[22,3] = 0.1038743148623453
$arr[22,4] = 0
$arr[22,5] = 0.002493689545182114
$arr[22,6] = 0
$arr[22,7] = 1.213889920155935
$arr[22,8] = 0
$arr[22,9] = 0.4226212513234486
$arr[22,10] = 0.2269310216957194
$arr[22,11] = 0.1775920323789038
$arr[22,12] = 0
$arr[22,13] = 4.946002961966542

$arr[23,0] = 0.537024338077515
$arr[23,1] = 0.1003495607393887
$arr[23,2] = 0.08824348106851687
$arr[23,3] = 0.1042584207893249
$arr[23,4] = 0
$arr[23,5] = 0.002498511920380854
$arr[23,6] = 0
$arr[23,7] = 1.222019079532402
$arr[23,8] = 0
$arr[23,9] = 0.357312121987519
$arr[23,10] = 0.2207937531578352
$arr[23,11] = 0.1644097249493193
$arr[23,12] = 0
$arr[23,13] = 4.95620917005229

$rng = $ws.Range($ws.Cells.Item(2, 2), $ws.Cells.Item(25, 15))
$rng.Value = $arr

Write-Output "Done updating pl_mw values"